$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("iOS_Estimate")

# 2a: "1. Implement the request body..." row (H13) status moves from "in progress" to "completed"
$ws.Range("H13").Value = "completed"

# 7a: new sub-task row for "Getting the details of parameters from WIFI..." (row 34)
$ws.Range("A34").Value = "a.created ipa to determine parameters of wifi"
$ws.Range("A34").IndentLevel = 6
$ws.Range("A34").HorizontalAlignment = -4131
$ws.Range("B34").Value = 3
$ws.Range("E34").Value = 3

# Update totals row (row 36) to reflect the added 3 dev hours / 3 total hours
$ws.Range("B36").Value = 59
$ws.Range("E36").Value = 73

# Move active selection to E37, matching the saved workbook view state
$ws.Range("E37").Select()
